$d = $word.ActiveDocument

# --- Change 1: insert sentence after "RDS' Service-Linked Admin Secrets." ---
$rng = $d.Content
$found = $rng.Find.Execute("RDS’ Service-Linked Admin Secrets", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
Write-Host "Change1 Found:" $found

# Walk forward character-by-character from the end of the hyperlink match to
# locate the literal "." that follows it (there may be 1+ hidden/zero-width
# boundary positions between the hyperlink run and the following run).
$scanPos = $rng.End
for ($i = 0; $i -lt 6; $i++) {
    $probe = $d.Range($scanPos, $scanPos + 1)
    if ($probe.Text -eq ".") {
        break
    }
    $scanPos = $scanPos + 1
}
$dotEnd = $scanPos + 1
Write-Host "Period located at:" $scanPos "-" $dotEnd

$afterDot = $d.Range($dotEnd, $dotEnd)
Write-Host "afterDot Text before:" "[$($afterDot.Text)]"
$afterDot.InsertBefore(" All new RDS instances since launch use this integration by default when rotation is configured on any database users.")

# --- Change 2: "Wrote an " -> "Published an " ---
$found2 = $d.Content.Find.Execute("Wrote an ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Published an ", 2)
Write-Host "Change2 Found:" $found2

# --- Change 3: "Developed a data key" -> "Implemented a data key" ---
$found3 = $d.Content.Find.Execute("Developed a data key", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Implemented a data key", 2)
Write-Host "Change3 Found:" $found3

# --- Change 4: "...and wrote a " -> "...and authored a " ---
$found4 = $d.Content.Find.Execute("Created a hybrid secret rotation strategy to meet compliance requirements and wrote a ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Created a hybrid secret rotation strategy to meet compliance requirements and authored a ", 2)
Write-Host "Change4 Found:" $found4

# --- Change 5: "Built and presented 10 lectures..." -> "Designed a curriculum with 10 lectures..." ---
$found5 = $d.Content.Find.Execute("Built and presented 10 lectures on Git, Code Reviews, Unit Testing, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Designed a curriculum with 10 lectures on Git, Code Reviews, Unit Testing, ", 2)
Write-Host "Change5 Found:" $found5

# --- Change 6: "The app also allows users to create an account..." -> "Users can also create an account..." ---
$found6 = $d.Content.Find.Execute("The app also allows users to create an account and save their pies for later.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Users can also create an account and save their pies for later.", 2)
Write-Host "Change6 Found:" $found6

# --- Change 7: ", Relative Valuation, etc.) " -> ", Relative Valuation) " ---
$found7 = $d.Content.Find.Execute(", Relative Valuation, etc.) ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", Relative Valuation) ", 2)
Write-Host "Change7 Found:" $found7
